$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Helper: locate a paragraph whose full (paragraph-mark-terminated)
# text exactly matches the given string, scanning the whole document.
# Using an exact-text scan (rather than a hard-coded paragraph index)
# keeps the script resilient to any incidental paragraph-count drift.
# ------------------------------------------------------------------
function Find-ParagraphByText($doc, [string]$target) {
    $paras = $doc.Paragraphs
    $count = $paras.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text -eq ($target + "`r")) {
            return $p
        }
    }
    return $null
}

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# ------------------------------------------------------------------
# Edit 1: "{{ etapa.titulo }}" paragraph.
# Drop the "Ttulo2" (Heading 2) paragraph style and instead apply a
# manual left alignment plus a single underline on the paragraph mark.
# ------------------------------------------------------------------
$pTitulo = Find-ParagraphByText $d "{{ etapa.titulo }}"
if ($pTitulo -ne $null) {
    $bodyTitulo = '<w:p w14:paraId="18CDBF64" w14:textId="77777777" w:rsidR="006A7BA9" w:rsidRPr="00841931" w:rsidRDefault="006A7BA9" w:rsidP="00841931">' + `
        '<w:pPr><w:jc w:val="left"/><w:rPr><w:u w:val="single"/></w:rPr></w:pPr>' + `
        '<w:proofErr w:type="gramStart"/>' + `
        '<w:r w:rsidRPr="00841931"><w:t xml:space="preserve">{{ </w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/>' + `
        '<w:r w:rsidRPr="00841931"><w:t>etapa</w:t></w:r>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
        '<w:r w:rsidRPr="00841931"><w:t>.</w:t></w:r>' + `
        '<w:proofErr w:type="gramStart"/>' + `
        '<w:r w:rsidRPr="00841931"><w:t>titulo</w:t></w:r>' + `
        '<w:proofErr w:type="spellEnd"/>' + `
        '<w:r w:rsidRPr="00841931"><w:t xml:space="preserve"> }</w:t></w:r>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
        '<w:r w:rsidRPr="00841931"><w:t>}</w:t></w:r>' + `
        '</w:p>'
    [void]$pTitulo.Range.InsertXML($pkgHeader + '<w:body>' + $bodyTitulo + '</w:body>' + $pkgFooter)
}

# ------------------------------------------------------------------
# Edit 2: "{% for item in etapa.itens %}" paragraph.
# Split the single run into three runs (wrapping "etapa.itens" with
# gramStart/gramEnd proofErr markers) while keeping the visible text
# and paragraph formatting identical.
# ------------------------------------------------------------------
$pItens = Find-ParagraphByText $d "{% for item in etapa.itens %}"
if ($pItens -ne $null) {
    $bodyItens = '<w:p w14:paraId="72D322B6" w14:textId="77777777" w:rsidR="006A7BA9" w:rsidRPr="00841931" w:rsidRDefault="006A7BA9" w:rsidP="006A7BA9">' + `
        '<w:pPr><w:jc w:val="left"/></w:pPr>' + `
        '<w:r w:rsidRPr="00841931"><w:t xml:space="preserve">{% for item in </w:t></w:r>' + `
        '<w:proofErr w:type="gramStart"/>' + `
        '<w:r w:rsidRPr="00841931"><w:t>etapa.itens</w:t></w:r>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
        '<w:r w:rsidRPr="00841931"><w:t xml:space="preserve"> %}</w:t></w:r>' + `
        '</w:p>'
    [void]$pItens.Range.InsertXML($pkgHeader + '<w:body>' + $bodyItens + '</w:body>' + $pkgFooter)
}

Write-Output "edit complete"
